{"js": "// Update the lattice-multiplication exercise table: every cell keeps its\n// existing layout (problem line, spaced-digit line, \"----\" rule line, two\n// carry rows) but the numbers themselves change. Each table cell is a\n// single run whose five lines are separated by manual line breaks\n// (<w:br/>), so we rebuild each cell's text as one string joined with the\n// vertical-tab character (\"\\u000B\") Office.js uses to represent a <w:br/>\n// inside Range/values text, and write it back with Range.insertText so the\n// run keeps its original formatting (sz=32) and Word re-adds\n// xml:space=\"preserve\" on the lines that need it.\n\n// New 5-line contents for every cell, addressed as [row][col] (0-based),\n// matching the table's current 5 rows x 3 columns.\nconst afterCells = [\n  [\n    [\"61 x 69\", \"  6    9\", \"  ----\", \"6|    |\", \"1|    |\"],\n    [\"60 x 41\", \"  4    1\", \"  ----\", \"6|    |\", \"0|    |\"],\n    [\"60 x 20\", \"  2    0\", \"  ----\", \"6|    |\", \"0|    |\"],\n  ],\n  [\n    [\"54 x 57\", \"  5    7\", \"  ----\", \"5|    |\", \"4|    |\"],\n    [\"81 x 76\", \"  7    6\", \"  ----\", \"8|    |\", \"1|    |\"],\n    [\"56 x 98\", \"  9    8\", \"  ----\", \"5|    |\", \"6|    |\"],\n  ],\n  [\n    [\"20 x 60\", \"  6    0\", \"  ----\", \"2|    |\", \"0|    |\"],\n    [\"74 x 56\", \"  5    6\", \"  ----\", \"7|    |\", \"4|    |\"],\n    [\"62 x 22\", \"  2    2\", \"  ----\", \"6|    |\", \"2|    |\"],\n  ],\n  [\n    [\"85 x 40\", \"  4    0\", \"  ----\", \"8|    |\", \"5|    |\"],\n    [\"39 x 53\", \"  5    3\", \"  ----\", \"3|    |\", \"9|    |\"],\n    [\"94 x 87\", \"  8    7\", \"  ----\", \"9|    |\", \"4|    |\"],\n  ],\n  [\n    [\"42 x 29\", \"  2    9\", \"  ----\", \"4|    |\", \"2|    |\"],\n    [\"57 x 80\", \"  8    0\", \"  ----\", \"5|    |\", \"7|    |\"],\n    [\"78 x 95\", \"  9    5\", \"  ----\", \"7|    |\", \"8|    |\"],\n  ],\n];\n\nconst LINE_BREAK = \"\\u000B\"; // Office.js's textual stand-in for <w:br/>\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nfor (let r = 0; r < afterCells.length; r++) {\n  for (let c = 0; c < afterCells[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(afterCells[r][c].join(LINE_BREAK), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication exercise table: every cell keeps its\n# existing layout (problem line, spaced-digit line, \"----\" rule line, two\n# carry rows) but the numbers themselves change. Each cell's text is set in\n# one shot via Range.Text, joining the five lines with a vertical-tab\n# character (chr(11)), which Word's COM layer turns back into manual line\n# breaks (<w:br/>) inside the cell's single run - preserving the run's\n# existing formatting (sz=32) and re-adding xml:space=\"preserve\" where the\n# text needs it.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$nl = [char]11\n\n# New 5-line contents for every cell, addressed as (row, col) with Word's\n# 1-based indices, matching the table's current 5 rows x 3 columns.\n$afterCells = @{\n  \"1,1\" = @(\"61 x 69\", \"  6    9\", \"  ----\", \"6|    |\", \"1|    |\")\n  \"1,2\" = @(\"60 x 41\", \"  4    1\", \"  ----\", \"6|    |\", \"0|    |\")\n  \"1,3\" = @(\"60 x 20\", \"  2    0\", \"  ----\", \"6|    |\", \"0|    |\")\n\n  \"2,1\" = @(\"54 x 57\", \"  5    7\", \"  ----\", \"5|    |\", \"4|    |\")\n  \"2,2\" = @(\"81 x 76\", \"  7    6\", \"  ----\", \"8|    |\", \"1|    |\")\n  \"2,3\" = @(\"56 x 98\", \"  9    8\", \"  ----\", \"5|    |\", \"6|    |\")\n\n  \"3,1\" = @(\"20 x 60\", \"  6    0\", \"  ----\", \"2|    |\", \"0|    |\")\n  \"3,2\" = @(\"74 x 56\", \"  5    6\", \"  ----\", \"7|    |\", \"4|    |\")\n  \"3,3\" = @(\"62 x 22\", \"  2    2\", \"  ----\", \"6|    |\", \"2|    |\")\n\n  \"4,1\" = @(\"85 x 40\", \"  4    0\", \"  ----\", \"8|    |\", \"5|    |\")\n  \"4,2\" = @(\"39 x 53\", \"  5    3\", \"  ----\", \"3|    |\", \"9|    |\")\n  \"4,3\" = @(\"94 x 87\", \"  8    7\", \"  ----\", \"9|    |\", \"4|    |\")\n\n  \"5,1\" = @(\"42 x 29\", \"  2    9\", \"  ----\", \"4|    |\", \"2|    |\")\n  \"5,2\" = @(\"57 x 80\", \"  8    0\", \"  ----\", \"5|    |\", \"7|    |\")\n  \"5,3\" = @(\"78 x 95\", \"  9    5\", \"  ----\", \"7|    |\", \"8|    |\")\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n  for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n    $lines = $afterCells[\"$r,$c\"]\n    $tbl.Cell($r, $c).Range.Text = [string]::Join($nl, $lines)\n  }\n}\n"}
